# Fix bug where interpolation crashed.
#
# The "carbon_intensity" per-country sheet actually held stale
# "power_latop"-style numeric overrides; it should be the "time_laptop"
# per-country overrides (interp JSON), and the params sheet had the
# "power_latop" / "carbon_intensity" rows' data swapped by mistake.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the "carbon_intensity" sheet -> "time_laptop"
# ------------------------------------------------------------------
$params = $wb.Worksheets.Item("params")
$timeLaptop = $wb.Worksheets.Item("carbon_intensity")
$timeLaptop.Name = "time_laptop"
$changes = $wb.Worksheets.Item("changes")
$metadata = $wb.Worksheets.Item("metadata")

# ------------------------------------------------------------------
# 2. params sheet: fix the swapped power_latop / carbon_intensity rows
#    row 2 was "power_latop" -> should be "carbon_intensity"
#    row 6 was "carbon_intensity" -> should be "power_latop"
# ------------------------------------------------------------------

# --- row 2 becomes carbon_intensity ---
$params.Range("A2").Value = "carbon_intensity"
$params.Range("C2").Value = "exp"
$params.Range("D2").Clear()
$params.Range("E2").Value = 0.5
$params.Range("F2").Value = -0.1
$params.Range("G2").Value = 0.1
$params.Range("H2").Value = 0.05
$params.Range("J2").Value = "kg/kWh"
$params.Range("P2").Clear()
$params.Range("Q2").Clear()
$params.Range("R2").Clear()
$params.Range("S2").Value = 5

# --- row 6 becomes power_latop ---
$params.Range("A6").Value = "power_latop"
$params.Range("C6").Value = "interp"
$params.Range("D6").Value = "linear"
$params.Range("E6").Value = '{"2020-01-01":10, "2031-06-01":9.5}'
$params.Range("F6").Value = 0
$params.Range("G6").Value = 4
$params.Range("H6").Value = 0.05
$params.Range("J6").Value = "W"
$params.Range("P6").Style = $params.Range("P2").Style
$params.Range("P6").Value = "what does it mean? How do collect this info?"
$params.Range("Q6").Style = $params.Range("Q2").Style
$params.Range("Q6").Value = "x"
$params.Range("R6").Style = $params.Range("R2").Style
$params.Range("R6").Value = "power draw of laptop"
$params.Range("S6").Value = 0

# ------------------------------------------------------------------
# 3. params sheet: drop the two trailing blank rows (8 & 9)
# ------------------------------------------------------------------
$params.Rows.Item(8).Delete()
$params.Rows.Item(8).Delete()

# ------------------------------------------------------------------
# 4. time_laptop (per-country) sheet: the UK/DE rows had plain numeric
#    "ref value" overrides left over from carbon_intensity; they need
#    the interp JSON (matching params!time_laptop) that was missing,
#    which is what caused interpolation to crash.
# ------------------------------------------------------------------
$timeLaptop.Range("A2").Value = "UK"
$timeLaptop.Range("C2").Value = '{"2020-01-01":100, "2031-06-01":95}'
$timeLaptop.Range("D2").Value = 0
$timeLaptop.Range("E2").Value = 5
$timeLaptop.Range("F2").Value = 0.05
$timeLaptop.Range("G2").NumberFormat = "0"
$timeLaptop.Range("G2").Value = 6

$timeLaptop.Range("A3").Value = "DE"
$timeLaptop.Range("C3").Value = '{"2020-01-01":100, "2031-06-01":95}'
$timeLaptop.Range("D3").Value = 1
$timeLaptop.Range("E3").Value = 5
$timeLaptop.Range("F3").Value = 0.05
$timeLaptop.Range("G3").Value = 7

# ------------------------------------------------------------------
# 5. Selections / active sheet bookkeeping
# ------------------------------------------------------------------
$params.Activate()
$params.Range("A6:S6").Select()

$timeLaptop.Range("D4").Select()

$changes.Range("F9").Select()

$params.Activate()
